$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 251.68
$ws.Range("I15").Value = 251.68
$ws.Range("K15").Value = 755.04
$ws.Range("M15").Value = -586.04

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 29646680
$ws.Range("I98").Value = 10527213
$ws.Range("J98").Value = 102300650
$ws.Range("K98").Value = 10527213
$ws.Range("L98").Value = 102300650
$ws.Range("M98").Value = -10525715
$ws.Range("N98").Value = -102303646

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 29646680
$ws.Range("I122").Value = 10527213
$ws.Range("J122").Value = 102300650
$ws.Range("K122").Value = 31581639
$ws.Range("L122").Value = 306901950
$ws.Range("M122").Value = -31579189
$ws.Range("N122").Value = -306906850

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2332972
$ws.Range("I125").Value = 800
$ws.Range("J125").Value = 2592102.2
$ws.Range("K125").Value = 7200
$ws.Range("L125").Value = 23328919.8
$ws.Range("M125").Value = -4740
$ws.Range("N125").Value = -23333839.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 438053.78
$ws.Range("I129").Value = 2368
$ws.Range("J129").Value = 1254964.6
$ws.Range("K129").Value = 7104
$ws.Range("L129").Value = 3764893.8
$ws.Range("M129").Value = -2104
$ws.Range("N129").Value = -3774893.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3729.31
$ws.Range("I138").Value = 1728.2354
$ws.Range("J138").Value = 4139.1685
$ws.Range("K138").Value = 5184.706200000001
$ws.Range("L138").Value = 12417.5055
$ws.Range("M138").Value = -44.70620000000054
$ws.Range("N138").Value = -22697.5055

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2620510.2
$ws.Range("I32").Value = 3409506.5
$ws.Range("J32").Value = 34355.5
$ws.Range("K32").Value = 3409506.5
$ws.Range("L32").Value = 34355.5
$ws.Range("M32").Value = -3409219.5
$ws.Range("N32").Value = -34929.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 7630.727
$ws.Range("I37").Value = 3987.5
$ws.Range("J37").Value = 17346
$ws.Range("K37").Value = 3987.5
$ws.Range("L37").Value = 17346
$ws.Range("M37").Value = -3714.5
$ws.Range("N37").Value = -17892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 21239262
$ws.Range("I132").Value = 27618596
$ws.Range("J132").Value = 7354827
$ws.Range("K132").Value = 82855788
$ws.Range("L132").Value = 22064481
$ws.Range("M132").Value = -82853258
$ws.Range("N132").Value = -22069541

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1871.48
$ws.Range("I86").Value = 1881.9381
$ws.Range("K86").Value = 1881.9381
$ws.Range("M86").Value = -758.9381000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1871.48
$ws.Range("I89").Value = 1881.9381
$ws.Range("K89").Value = 9409.690500000001
$ws.Range("M89").Value = -3793.690500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1891.0588
$ws.Range("I94").Value = 1218.3077
$ws.Range("J94").Value = 4077.5
$ws.Range("K94").Value = 1218.3077
$ws.Range("L94").Value = 4077.5
$ws.Range("M94").Value = -767.3077000000001
$ws.Range("N94").Value = -4979.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17057.143
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 17057.143
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 17057.143
$ws.Range("N31").Value = -17647.143
$ws.Range("M31").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 17057.143
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 17057.143
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 17057.143
$ws.Range("N34").Value = -17461.143
$ws.Range("M34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 27187.375
$ws.Range("J51").Value = 27187.375
$ws.Range("L51").Value = 27187.375
$ws.Range("N51").Value = -28659.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 23333.334
$ws.Range("J59").Value = 23333.334
$ws.Range("L59").Value = 23333.334
$ws.Range("N59").Value = -25623.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18331.666
$ws.Range("J60").Value = 18331.666
$ws.Range("L60").Value = 18331.666
$ws.Range("N60").Value = -19353.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 27187.375
$ws.Range("J61").Value = 27187.375
$ws.Range("L61").Value = 27187.375
$ws.Range("N61").Value = -27883.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 19999
$ws.Range("J68").Value = 19999
$ws.Range("L68").Value = 19999
$ws.Range("N68").Value = -21497

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 19999
$ws.Range("J71").Value = 19999
$ws.Range("L71").Value = 59997
$ws.Range("N71").Value = -67485

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 26870.715
$ws.Range("J74").Value = 26870.715
$ws.Range("L74").Value = 26870.715
$ws.Range("N74").Value = -28618.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 26870.715
$ws.Range("J77").Value = 26870.715
$ws.Range("L77").Value = 80612.145
$ws.Range("N77").Value = -89348.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 911271.25
$ws.Range("I134").Value = 1550.7941
$ws.Range("J134").Value = 4004320.8
$ws.Range("K134").Value = 4652.3823
$ws.Range("L134").Value = 12012962.4
$ws.Range("M134").Value = -2117.3823
$ws.Range("N134").Value = -12018032.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2325119.5
$ws.Range("I5").Value = 5128803.5
$ws.Range("J5").Value = 1390558.1
$ws.Range("K5").Value = 15386410.5
$ws.Range("L5").Value = 4171674.3
$ws.Range("M5").Value = -15386298.5
$ws.Range("N5").Value = -4171898.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2710.2727
$ws.Range("I113").Value = 1801.9
$ws.Range("J113").Value = 3467.25
$ws.Range("K113").Value = 5405.700000000001
$ws.Range("L113").Value = 10401.75
$ws.Range("M113").Value = -3235.700000000001
$ws.Range("N113").Value = -14741.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 298.92682
$ws.Range("I122").Value = 302.84616
$ws.Range("J122").Value = 222.5
$ws.Range("K122").Value = 2725.61544
$ws.Range("L122").Value = 2002.5
$ws.Range("M122").Value = -275.61544
$ws.Range("N122").Value = -6902.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2325119.5
$ws.Range("I135").Value = 5128803.5
$ws.Range("J135").Value = 1390558.1
$ws.Range("K135").Value = 46159231.5
$ws.Range("L135").Value = 12515022.9
$ws.Range("M135").Value = -46156696.5
$ws.Range("N135").Value = -12520092.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 28822
$ws.Range("I122").Value = 41470.31
$ws.Range("J122").Value = 5332.2856
$ws.Range("K122").Value = 124410.93
$ws.Range("L122").Value = 15996.8568
$ws.Range("M122").Value = -121960.93
$ws.Range("N122").Value = -20896.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 15217
$ws.Range("J74").Value = 15217
$ws.Range("L74").Value = 15217
$ws.Range("N74").Value = -17213

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 15217
$ws.Range("J77").Value = 15217
$ws.Range("L77").Value = 45651
$ws.Range("N77").Value = -55635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3044196.5
$ws.Range("I132").Value = 3972381.8
$ws.Range("J132").Value = 6499.909
$ws.Range("K132").Value = 11917145.4
$ws.Range("L132").Value = 19499.727
$ws.Range("M132").Value = -11914615.4
$ws.Range("N132").Value = -24559.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 160.8
$ws.Range("I113").Value = 102
$ws.Range("K113").Value = 306
$ws.Range("M113").Value = 1864
